$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("_name"), shifting everything
# from B onward one column to the right (B:P -> C:Q).
$ws.Range("B1:B3").EntireColumn.Insert()

# Match the new column B's width to column A's (Excel's "insert column"
# carries over the left-neighbour's width).
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# Populate the new "_itemType" column.
$ws.Range("B1").Value2 = "_itemType"
$ws.Range("B2").Value2 = "Consumable"
$ws.Range("B3").Value2 = "Consumable"

# Data fix that rode along in the same commit: Power Potion's
# _stageDuration (now column Q) becomes 1.
$ws.Range("Q3").Value2 = 1

# Restore the active selection to match the saved view state.
$ws.Range("M10").Select() | Out-Null
